# Generate Report for Handback
# Updates the "Ready for handoff" row (321f485f-...) on both the zh-cn and
# de-de sheets with the freshly generated handback report info:
#   - Latest Target File (col I) becomes a hyperlink to the .md file
#   - Latest Handback File (col J) is filled in with the handback xlf name
#   - Latest Handback DateTime (col K) is filled in
#   - Error Detail (col P) explains the version mismatch
# Also widens columns I, J and P to match the other wide (40) columns.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bb1563339901e18b51bb252728e9ba7d50e9cf63/e2e/321f485f-333c-4625-8624-0a76dfaa6197.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ddd3043751cc03593ad28eb8aa404e6be9e5752f/e2e/321f485f-333c-4625-8624-0a76dfaa6197.md."
$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ddd3043751cc03593ad28eb8aa404e6be9e5752f/e2e/321f485f-333c-4625-8624-0a76dfaa6197.md"
$mdName = "321f485f-333c-4625-8624-0a76dfaa6197.md"

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("I5").Value = $mdName
$wsZh.Hyperlinks.Add($wsZh.Range("I5"), $targetUrl, "", "", $mdName) | Out-Null
$wsZh.Range("J5").Value = "321f485f-333c-4625-8624-0a76dfaa6197.e6eb9a2f14dbf57b524e83b0678c8dcc2167a5de.zh-cn.xlf"
$wsZh.Range("K5").Value = "2016-10-27 09:10:35"
$wsZh.Range("P5").Value = $errorDetail

$wsZh.Columns.Item(9).ColumnWidth = 39.17
$wsZh.Columns.Item(10).ColumnWidth = 39.17
$wsZh.Columns.Item(16).ColumnWidth = 39.17

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I5").Value = $mdName
$wsDe.Hyperlinks.Add($wsDe.Range("I5"), $targetUrl, "", "", $mdName) | Out-Null
$wsDe.Range("J5").Value = "321f485f-333c-4625-8624-0a76dfaa6197.e6eb9a2f14dbf57b524e83b0678c8dcc2167a5de.de-de.xlf"
$wsDe.Range("K5").Value = "2016-10-27 09:10:51"
$wsDe.Range("P5").Value = $errorDetail

$wsDe.Columns.Item(9).ColumnWidth = 39.17
$wsDe.Columns.Item(10).ColumnWidth = 39.17
$wsDe.Columns.Item(16).ColumnWidth = 39.17
